$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.479.21"
$ws.Range("E2").Value = "  +4.77%  "

$ws.Range("D3").Value = "3.591.55"
$ws.Range("E3").Value = "  +16.89%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'592.30"
$ws.Range("E5").Value = "  +3.26%  "

$ws.Range("D6").Value = "'184.94"
$ws.Range("E6").Value = "  +8.65%  "

$ws.Range("D7").Value = "3.579.54"
$ws.Range("E7").Value = "  +16.57%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  +4.54%  "

$ws.Range("E10").Value = "  +6.79%  "

$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "  +3.24%  "

$ws.Range("D12").Value = "'0.494"
$ws.Range("E12").Value = "  +5.61%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'38.72"
$ws.Range("E13").Value = "  +8.59%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000251"
$ws.Range("E14").Value = "  +5.32%  "

$ws.Range("D15").Value = "4.175.83"
$ws.Range("E15").Value = "  +16.70%  "

$ws.Range("D16").Value = "3.580.58"

$ws.Range("D17").Value = "69.540.78"
$ws.Range("E17").Value = "  +5.10%  "

$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +8.03%  "

$ws.Range("D20").Value = "'17.09"
$ws.Range("E20").Value = "  +3.03%  "

$ws.Range("D21").Value = "'506.45"
$ws.Range("E21").Value = "  +4.30%  "

$ws.Range("D22").Value = "'9.44"
$ws.Range("E22").Value = "  +23.37%  "

$ws.Range("D23").Value = "'0.742"
$ws.Range("E23").Value = "  +8.40%  "

$ws.Range("D24").Value = "'87.06"
$ws.Range("E24").Value = "  +5.76%  "

$ws.Range("D25").Value = "'13.44"
$ws.Range("E25").Value = "  +6.34%  "

$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +8.45%  "

$ws.Range("D27").Value = "'10.85"
$ws.Range("E27").Value = "  +6.56%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").Value = "'2.53"
$ws.Range("E29").Value = "  +12.45%  "

$ws.Range("D30").Value = "'8.12"
$ws.Range("E30").Value = "  +3.22%  "

$ws.Range("D31").Value = "'31.96"
$ws.Range("E31").Value = "  +15.42%  "

$ws.Range("E32").Value = "  +20.80%  "

$ws.Range("D33").Value = "'2.74"
$ws.Range("E33").Value = "  +5.75%  "

$ws.Range("D34").Value = "'0.117"
$ws.Range("E34").Value = "  +6.05%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'6.14"
$ws.Range("E36").Value = "  +10.68%  "

$ws.Range("E37").Value = "  +8.86%  "

$ws.Range("D38").Value = "'0.334"
$ws.Range("E38").Value = "  +11.87%  "

$ws.Range("D39").Value = "'46.90"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("D40").Value = "'2.10"
$ws.Range("E40").Value = "  +7.61%  "

$ws.Range("D41").Value = "'50.67"

$ws.Range("D42").Value = "'0.128"
$ws.Range("E42").Value = "  +4.64%  "

$ws.Range("D43").Value = "'8.78"
$ws.Range("E43").Value = "  +6.80%  "

$ws.Range("D44").Value = "3.077.80"
$ws.Range("E44").Value = "  +11.32%  "

$ws.Range("D45").Value = "'2.82"
$ws.Range("E45").Value = "  +11.24%  "

$ws.Range("D46").Value = "'404.44"
$ws.Range("E46").Value = "  +11.55%  "

$ws.Range("D47").Value = "'0.0364"
$ws.Range("E47").Value = "  +6.33%  "

$ws.Range("D48").Value = "'27.79"
$ws.Range("E48").Value = "  +14.67%  "

$ws.Range("D49").Value = "'135.08"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "'2.45"
$ws.Range("E51").Value = "  +14.68%  "

# Reset style on quote-prefixed numeric-looking text cells so they keep the default (unstyled) cell format
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
